{"js": "// 1. Make the title paragraph (\"Capstone 1 H1B LCA Petition Data Wrangling\") bold.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.font.bold = true;\n\n// 2. Remove the empty paragraph that immediately follows the title.\nconst emptyPara = paragraphs.items[1];\nemptyPara.delete();\n\nawait context.sync();\n\n// 3. Relocate the \"_GoBack\" bookmark from the end of the document to the\n//    point in the \"Other cleanup included...\" paragraph where the author's\n//    cursor last was (right in the middle of the word \"float\").\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst searchResults = body.search(\"flo\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst splitPoint = searchResults.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# 1. Make the title paragraph (\"Capstone 1 H1B LCA Petition Data Wrangling\") bold.\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.Font.Bold = 1\n\n# 2. Remove the empty paragraph that immediately follows the title.\n$emptyPara = $d.Paragraphs.Item(2)\n$emptyPara.Range.Delete()\n\n# 3. Relocate the \"_GoBack\" bookmark from the end of the document to the\n#    point in the \"Other cleanup included...\" paragraph where the author's\n#    cursor last was (right in the middle of the word \"float\").\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n$splitRange = $d.Content\n$splitRange.Find.Execute(\"flo\") | Out-Null\n$splitRange.Collapse(0)  # wdCollapseEnd\n\n$d.Bookmarks.Add(\"_GoBack\", $splitRange) | Out-Null\n"}
